$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(317, 8).WrapText = $true
$ws.Cells.Item(318, 8).Font.Bold = $true
$ws.Cells.Item(319, 8).Font.Italic = $true
